$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. "Handoff transform failed" -> "Ready for handoff"
#    This shared string is used on the Overview sheet (B2 & C2) and on
#    each language sheet's B2 ("Status" column for the .md source row).
# ------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("B2").Value = "Ready for handoff"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("B2").Value = "Ready for handoff"

# ------------------------------------------------------------------
# 2. zh-cn sheet: report the handed-off xlf file, its handoff time, and
#    flip the dependency row's handoff reason from "Ignored" to
#    "Include".
# ------------------------------------------------------------------
$wsZhCn.Range("C2").Value = "94524e43-f7a7-40fd-99f5-e98cc09d765c.dc6d88f9f785dd6c77b93aaa9a8c170e0523426f.zh-cn.xlf"
$wsZhCn.Range("C2").Style = "HyperLink"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("C2"), "https://github.com/OpenLocalizationTest/oltest/blob/d329c184eed499fd4b0db68276f0485cfcb7aed6/e2e/94524e43-f7a7-40fd-99f5-e98cc09d765c.dc6d88f9f785dd6c77b93aaa9a8c170e0523426f.zh-cn.xlf", "", "", "94524e43-f7a7-40fd-99f5-e98cc09d765c.dc6d88f9f785dd6c77b93aaa9a8c170e0523426f.zh-cn.xlf")

$wsZhCn.Range("D2").Value = "2016-01-28 09:35:20"
$wsZhCn.Range("H2").Value = "Include"

# ------------------------------------------------------------------
# 3. de-de sheet: same shape of update, with the de-de xlf file and its
#    own handoff time.
# ------------------------------------------------------------------
$wsDeDe.Range("C2").Value = "94524e43-f7a7-40fd-99f5-e98cc09d765c.dc6d88f9f785dd6c77b93aaa9a8c170e0523426f.de-de.xlf"
$wsDeDe.Range("C2").Style = "HyperLink"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("C2"), "https://github.com/OpenLocalizationTest/oltest/blob/d329c184eed499fd4b0db68276f0485cfcb7aed6/e2e/94524e43-f7a7-40fd-99f5-e98cc09d765c.dc6d88f9f785dd6c77b93aaa9a8c170e0523426f.de-de.xlf", "", "", "94524e43-f7a7-40fd-99f5-e98cc09d765c.dc6d88f9f785dd6c77b93aaa9a8c170e0523426f.de-de.xlf")

$wsDeDe.Range("D2").Value = "2016-01-28 09:35:33"
$wsDeDe.Range("H2").Value = "Include"
